$wb = $excel.ActiveWorkbook

# --- Summary ---
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B4").Value = "inf"
$ws.Range("B7").Value = 6889165.230378538
$ws.Range("B8").Value = 31295665.63009646
$ws.Range("B10").Value = 1007837.145333632

# --- Fed-in Capacity ---
$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Range("M2").Value = 0
$ws.Range("K3").Value = 168.166555548718
$ws.Range("N3").Value = 160.2368887416667
$ws.Range("O3").Value = 173.9674182222222
$ws.Range("P3").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("P6").Value = 163.4487770454829
$ws.Range("M8").Value = 281.0224045372727
$ws.Range("N8").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("M11").Value = 281.0224045372727
$ws.Range("M12").Value = 173.4035213848624
$ws.Range("O12").Value = 173.9674182222222
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = 160.2368887416667
$ws.Range("O15").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 173.4035213848624
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 163.4487770454829
$ws.Range("N20").Value = 0
$ws.Range("O21").Value = 173.9674182222222
$ws.Range("P21").Value = 163.4487770454829
$ws.Range("M23").Value = 281.0224045372727
$ws.Range("K24").Value = 0
$ws.Range("P24").Value = 163.4487770454829
$ws.Range("L27").Value = 0
$ws.Range("P27").Value = 0
$ws.Range("M29").Value = 281.0224045372727
$ws.Range("N29").Value = 279.8839375878409
$ws.Range("M30").Value = 173.4035213848624
$ws.Range("N30").Value = 0
$ws.Range("P30").Value = 163.4487770454829
$ws.Range("M32").Value = 0
$ws.Range("L33").Value = 169.0363433314465
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = 163.4487770454829
$ws.Range("K36").Value = 168.166555548718
$ws.Range("M36").Value = 173.4035213848624
$ws.Range("M38").Value = 0
$ws.Range("K39").Value = 168.166555548718
$ws.Range("N39").Value = 0
$ws.Range("P39").Value = 0
$ws.Range("N41").Value = 0
$ws.Range("L42").Value = 169.0363433314465
$ws.Range("M42").Value = 0
$ws.Range("O42").Value = 0
$ws.Range("P42").Value = 0
$ws.Range("O45").Value = 0
$ws.Range("P45").Value = 0

# --- Unmet Demand ---
$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Range("M2").Value = 281.0224045372727
$ws.Range("K3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 163.4487770454829
$ws.Range("N6").Value = 160.2368887416667
$ws.Range("P6").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 279.8839375878409
$ws.Range("L9").Value = 169.0363433314465
$ws.Range("P9").Value = 163.4487770454829
$ws.Range("M11").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("L15").Value = 169.0363433314465
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 173.9674182222222
$ws.Range("K18").Value = 168.166555548718
$ws.Range("L18").Value = 169.0363433314465
$ws.Range("M18").Value = 0
$ws.Range("O18").Value = 173.9674182222222
$ws.Range("P18").Value = 0
$ws.Range("N20").Value = 279.8839375878409
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("M23").Value = 0
$ws.Range("K24").Value = 168.166555548718
$ws.Range("P24").Value = 0
$ws.Range("L27").Value = 169.0363433314465
$ws.Range("P27").Value = 163.4487770454829
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = 160.2368887416667
$ws.Range("P30").Value = 0
$ws.Range("M32").Value = 281.0224045372727
$ws.Range("L33").Value = 0
$ws.Range("O33").Value = 173.9674182222222
$ws.Range("P33").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").Value = 0
$ws.Range("M38").Value = 281.0224045372727
$ws.Range("K39").Value = 0
$ws.Range("N39").Value = 160.2368887416667
$ws.Range("P39").Value = 163.4487770454829
$ws.Range("N41").Value = 279.8839375878409
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 173.4035213848624
$ws.Range("O42").Value = 173.9674182222222
$ws.Range("P42").Value = 163.4487770454829
$ws.Range("O45").Value = 173.9674182222222
$ws.Range("P45").Value = 163.4487770454829

# --- Household Surplus ---
$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B2").Value = 76263.75050442164
$ws.Range("B3").Value = 76460.23147123269
$ws.Range("B4").Value = 77399.30309260261
$ws.Range("B5").Value = 152947.1086751887
$ws.Range("B6").Value = 127496.2295679065
$ws.Range("B7").Value = 63982.27304951782
$ws.Range("B8").Value = 76725.1993160991
$ws.Range("B9").Value = 51274.32020881689
$ws.Range("B10").Value = 90023.83060265993
$ws.Range("B11").Value = 128149.1490984682
$ws.Range("B12").Value = 89002.72814642602
$ws.Range("B13").Value = 100892.3052910578
$ws.Range("B14").Value = 51465.79092770109
$ws.Range("B15").Value = 64166.87604895027
$ws.Range("B16").Value = 64396.87046793634
